# Quarterly indexing esoteric bug-fix operation
#
# Column A holds dates that mark the 1st day of a calendar quarter
# (Jan/Apr/Jul/Oct 1st). The fix re-indexes each date to the 15th of the
# following month (i.e. mid-quarter) instead of the 1st of the quarter's
# first month. Column B (the GDP index values) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 150; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldDate = $cell.Value()
    $shifted = $oldDate.AddMonths(1)
    $newDate = Get-Date -Year $shifted.Year -Month $shifted.Month -Day 15 -Hour 0 -Minute 0 -Second 0
    $cell.Value = $newDate
}
